$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1132.8182
$ws.Range("J17").Value = 1132.8182
$ws.Range("L17").Value = 3398.4546
$ws.Range("N17").Value = -3734.4546

$ws.Range("H100").Value = 2872.1428
$ws.Range("I100").Value = 2621
$ws.Range("J100").Value = 3500
$ws.Range("K100").Value = 2621
$ws.Range("L100").Value = 3500
$ws.Range("M100").Value = -2080
$ws.Range("N100").Value = -4582

$ws.Range("H121").Value = 773.2
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 773.2
$ws.Range("K121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("M121").Value = 2319.6
$ws.Range("N121").Value = -5813.6

$ws.Range("H138").Value = 2859593
$ws.Range("I138").Value = 2181
$ws.Range("J138").Value = 3776121.2
$ws.Range("K138").Value = 6543
$ws.Range("L138").Value = 11328363.6
$ws.Range("M138").Value = -1403
$ws.Range("N138").Value = -11338643.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2215.2593
$ws.Range("I2").Value = 1912.4706
$ws.Range("J2").Value = 2730
$ws.Range("K2").Value = 1912.4706
$ws.Range("L2").Value = 2730
$ws.Range("M2").Value = -1799.4706
$ws.Range("N2").Value = -2956

$ws.Range("H15").Value = 1100
$ws.Range("J15").Value = 1100
$ws.Range("L15").Value = 1100
$ws.Range("N15").Value = -1800

$ws.Range("H32").Value = 7730.49
$ws.Range("I32").Value = 5992.388
$ws.Range("J32").Value = 17579.732
$ws.Range("K32").Value = 5992.388
$ws.Range("L32").Value = 17579.732
$ws.Range("M32").Value = -5705.388
$ws.Range("N32").Value = -18153.732

$ws.Range("H61").Value = 38540116
$ws.Range("I61").Value = 52685356
$ws.Range("J61").Value = 145885.72
$ws.Range("K61").Value = 52685356
$ws.Range("L61").Value = 145885.72
$ws.Range("M61").Value = -52685144
$ws.Range("N61").Value = -146309.72

$ws.Range("H70").Value = 39998
$ws.Range("J70").Value = 39998
$ws.Range("L70").Value = 39998
$ws.Range("N70").Value = -40538

$ws.Range("H73").Value = 39998
$ws.Range("J73").Value = 39998
$ws.Range("L73").Value = 39998
$ws.Range("N73").Value = -41870

$ws.Range("H97").Value = 2842156.8
$ws.Range("I97").Value = 4167869.2
$ws.Range("J97").Value = 1344.4286
$ws.Range("K97").Value = 4167869.2
$ws.Range("L97").Value = 1344.4286
$ws.Range("M97").Value = -4167373.2
$ws.Range("N97").Value = -2336.4286

$ws.Range("H102").Value = 7938712
$ws.Range("I102").Value = 10206233
$ws.Range("J102").Value = 2388
$ws.Range("K102").Value = 10206233
$ws.Range("L102").Value = 2388
$ws.Range("M102").Value = -10204611
$ws.Range("N102").Value = -5632

$ws.Range("H110").Value = 1516.3684
$ws.Range("I110").Value = 1073
$ws.Range("J110").Value = 2009
$ws.Range("K110").Value = 1073
$ws.Range("L110").Value = 2009
$ws.Range("M110").Value = 972
$ws.Range("N110").Value = -6099

$ws.Range("H116").Value = 2215.2593
$ws.Range("I116").Value = 1912.4706
$ws.Range("J116").Value = 2730
$ws.Range("K116").Value = 1912.4706
$ws.Range("L116").Value = 2730
$ws.Range("M116").Value = 381.5293999999999
$ws.Range("N116").Value = -7318

$ws.Range("H132").Value = 21833040
$ws.Range("I132").Value = 29473714
$ws.Range("J132").Value = 184466.5
$ws.Range("K132").Value = 88421142
$ws.Range("L132").Value = 553399.5
$ws.Range("M132").Value = -88418612
$ws.Range("N132").Value = -558459.5

$ws.Range("H136").Value = 38540116
$ws.Range("I136").Value = 52685356
$ws.Range("J136").Value = 145885.72
$ws.Range("K136").Value = 158056068
$ws.Range("L136").Value = 437657.16
$ws.Range("M136").Value = -158053518
$ws.Range("N136").Value = -442757.16

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2215.2593
$ws.Range("I3").Value = 1912.4706
$ws.Range("J3").Value = 2730
$ws.Range("K3").Value = 1912.4706
$ws.Range("L3").Value = 2730
$ws.Range("M3").Value = -1798.4706
$ws.Range("N3").Value = -2958

$ws.Range("H134").Value = 4904138.5
$ws.Range("I134").Value = 2173.5715
$ws.Range("J134").Value = 27779974
$ws.Range("K134").Value = 6520.7145
$ws.Range("L134").Value = 83339922
$ws.Range("M134").Value = -3985.7145
$ws.Range("N134").Value = -83344992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 374.75
$ws.Range("I19").Value = 374.75
$ws.Range("K19").Value = 374.75
$ws.Range("M19").Value = -204.75

$ws.Range("H24").Value = 374.75
$ws.Range("I24").Value = 374.75
$ws.Range("K24").Value = 374.75
$ws.Range("M24").Value = -204.75

$ws.Range("H132").Value = 40362.77
$ws.Range("I132").Value = 1507.9412
$ws.Range("J132").Value = 113755.22
$ws.Range("K132").Value = 4523.8236
$ws.Range("L132").Value = 341265.66
$ws.Range("M132").Value = -1993.8236
$ws.Range("N132").Value = -346325.66

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 100
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 300
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -127

$ws.Range("H75").Value = 4137.1816
$ws.Range("I75").Value = 3332
$ws.Range("J75").Value = 4217.7
$ws.Range("K75").Value = 9996
$ws.Range("L75").Value = 12653.1
$ws.Range("M75").Value = -8998
$ws.Range("N75").Value = -14649.1

$ws.Range("H76").Value = 3407.1428
$ws.Range("I76").Value = 1250
$ws.Range("J76").Value = 3766.6667
$ws.Range("K76").Value = 3750
$ws.Range("L76").Value = 11300.0001
$ws.Range("M76").Value = -3367
$ws.Range("N76").Value = -12066.0001

$ws.Range("H78").Value = 4137.1816
$ws.Range("I78").Value = 3332
$ws.Range("J78").Value = 4217.7
$ws.Range("K78").Value = 29988
$ws.Range("L78").Value = 37959.3
$ws.Range("M78").Value = -24996
$ws.Range("N78").Value = -47943.3

$ws.Range("H79").Value = 3407.1428
$ws.Range("I79").Value = 1250
$ws.Range("J79").Value = 3766.6667
$ws.Range("K79").Value = 3750
$ws.Range("L79").Value = 11300.0001
$ws.Range("M79").Value = -2424
$ws.Range("N79").Value = -13952.0001

$ws.Range("H88").Value = 2968.75
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 2968.75
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = 8906.25
$ws.Range("N88").Value = -9762.25

$ws.Range("H91").Value = 2968.75
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 2968.75
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = 8906.25
$ws.Range("N91").Value = -11870.25

$ws.Range("H129").Value = 5210476
$ws.Range("I129").Value = 1309.875
$ws.Range("J129").Value = 10419642
$ws.Range("K129").Value = 3929.625
$ws.Range("L129").Value = 31258926
$ws.Range("M129").Value = 1070.375
$ws.Range("N129").Value = -31268926

$ws.Range("H141").Value = 11892.95
$ws.Range("I141").Value = 2228.7778
$ws.Range("J141").Value = 19800
$ws.Range("K141").Value = 6686.3334
$ws.Range("L141").Value = 59400
$ws.Range("M141").Value = -1506.3334
$ws.Range("N141").Value = -69760

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2793.4
$ws.Range("I40").Value = 2741.75
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2741.75
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -2605.75
$ws.Range("N40").Value = -3272

$ws.Range("H122").Value = 3673.9119
$ws.Range("I122").Value = 4318.0713
$ws.Range("J122").Value = 3223
$ws.Range("K122").Value = 12954.2139
$ws.Range("L122").Value = 9669
$ws.Range("M122").Value = -10504.2139
$ws.Range("N122").Value = -14569

$ws.Range("H132").Value = 185709.81
$ws.Range("I132").Value = 4726
$ws.Range("J132").Value = 668333.3
$ws.Range("K132").Value = 14178
$ws.Range("L132").Value = 2004999.9
$ws.Range("M132").Value = -11648
$ws.Range("N132").Value = -2010059.9

$ws.Range("H137").Value = 30600
$ws.Range("I137").Value = 23400
$ws.Range("J137").Value = 32200
$ws.Range("K137").Value = 23400
$ws.Range("L137").Value = 32200
$ws.Range("M137").Value = -18300
$ws.Range("N137").Value = -42400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 3602.4
$ws.Range("I26").Value = 2337.3333
$ws.Range("J26").Value = 5500
$ws.Range("K26").Value = 2337.3333
$ws.Range("L26").Value = 5500
$ws.Range("M26").Value = -2044.3333
$ws.Range("N26").Value = -6086

$ws.Range("H81").Value = 1937.875
$ws.Range("I81").Value = 1071.4286
$ws.Range("J81").Value = 2611.7778
$ws.Range("K81").Value = 2142.8572
$ws.Range("L81").Value = 5223.5556
$ws.Range("M81").Value = -1081.8572
$ws.Range("N81").Value = -7345.5556

$ws.Range("H84").Value = 1937.875
$ws.Range("I84").Value = 1071.4286
$ws.Range("J84").Value = 2611.7778
$ws.Range("K84").Value = 10714.286
$ws.Range("L84").Value = 26117.778
$ws.Range("M84").Value = -5410.286
$ws.Range("N84").Value = -36725.778

$ws.Range("H132").Value = 75651.44500000001
$ws.Range("I132").Value = 60146.176
$ws.Range("J132").Value = 102010.4
$ws.Range("K132").Value = 180438.528
$ws.Range("L132").Value = 306031.2
$ws.Range("M132").Value = -177908.528
$ws.Range("N132").Value = -311091.2

$ws.Range("H136").Value = 45011.53
$ws.Range("I136").Value = 34837.566
$ws.Range("J136").Value = 62965.59
$ws.Range("K136").Value = 104512.698
$ws.Range("L136").Value = 188896.77
$ws.Range("M136").Value = -101962.698
$ws.Range("N136").Value = -193996.77
